$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Ev. Informe. En relación", "Evacua Informe"),
    @("Inadmisible", "Incompetencia"),
    @("Traslado", "Ica Informa"),
    @("Agréguese Extraordinariamente a la tabla del 03 de marzo", "Agréguese A Tabla"),
    @("Ev. Informe. Pase al sr. presidente", "Evacua Informe"),
    @("Se declara admisible, 05 días hábiles. Concede ONI.", "Admisibles"),
    @("Nota de Acuerdo: redacción Margarita Sanhueza", "Acuerdo"),
    @("Concede ampliación de plazo por el término de 05 días", "Ica Amplia Plazo"),
    @("Concede ampliación de plazo por el término de 08 días", "Ica Amplia Plazo"),
    @("Pide informe cumplimiento de sentencia a la recurrida", "Ica Solicita Diligencia"),
    @("Concede ampliación de Plazo por el término de 05 días", "Ica Amplia Plazo"),
    @("No ha lugar ocúrrase ante quien corresponda", "Ica Informa"),
    @("Se declara admisible, 10 días. NHL ONI.", "Admisibles"),
    @("Concede recurso, interconexión", "Concede Apelación")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
